$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    'D2' = '256.98'
    'E2' = '0.04%'
    'D3' = '26.73'
    'E3' = '-1.17%'
    'D4' = '4.636'
    'E4' = '0.34%'
    'D5' = '0.05935'
    'E5' = '0.67%'
    'E6' = '-0.40%'
    'D7' = '0.8567'
    'E7' = '-0.96%'
    'D8' = '0.9163'
    'E8' = '-2.20%'
    'D9' = '0.1382'
    'E9' = '-1.63%'
    'D10' = '0.04445'
    'E10' = '15.61%'
    'D11' = '0.06997'
    'E11' = '-1.05%'
    'D12' = '0.03012'
    'E12' = '-5.97%'
    'D13' = '0.09108'
    'E13' = '-1.55%'
    'D14' = '0.001536'
    'E14' = '0.01%'
    'D15' = '0.0006019'
    'E15' = '0.24%'
    'D16' = '0.006162'
    'E16' = '0.71%'
    'E17' = '-1.37%'
    'D18' = '3.131'
    'E18' = '-1.82%'
    'D20' = '0.3095'
    'E20' = '-0.07%'
    'E21' = '1.61%'
    'D22' = '3.864'
    'E22' = '0.25%'
    'D23' = '0.04185'
    'E23' = '-0.93%'
    'D24' = '0.001215'
    'E24' = '-0.20%'
    'D25' = '0.004476'
    'E25' = '4.63%'
    'D26' = '0.0001198'
    'E26' = '-0.01%'
    'D27' = '0.0001713'
    'E27' = '-11.49%'
    'D40' = '0.03814'
    'E40' = '-0.38%'
    'D41' = '0.1101'
    'E41' = '0.15%'
    'D42' = '0.003702'
    'E42' = '-39.63%'
    'D43' = '0.002426'
    'E43' = '4.99%'
    'D44' = '0.01507'
    'E44' = '29.73%'
    'D45' = '0.00005093'
    'E45' = '-6.62%'
    'D46' = '0.00000000749'
    'E46' = '0.01%'
    'D47' = '0.05004'
    'E47' = '-35.63%'
    'E48' = '10,479.89%'
    'D49' = '0.00002098'
    'E49' = '0.01%'
    'D50' = '0.0001998'
    'E50' = '0.01%'
}

foreach ($ref in $changes.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$ref]
    $cell.ClearFormats()
}
